{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block (and the\n// blank paragraph immediately preceding it) that followed the last\n// bibliography entry (\"SERAFINI, Maria Jos\u00e9...\").\n//\n// The document ends with:\n//   ... SERAFINI, Maria Jos\u00e9. Como escrever textos. ...\n//   (blank paragraph)\n//   Ver no Jupiter Salvar em pdf Salvar em docx\n//   \u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ...\n//   (blank paragraph)\n//   (paragraph with page-break-before)\n//\n// After the edit the two footer paragraphs and the blank paragraph that\n// introduced them are gone, so the bibliography's last entry is followed\n// directly by the blank paragraph + page-break paragraph that used to close\n// the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet lastBiblioIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"SERAFINI\") !== -1 && t.indexOf(\"Como escrever textos\") !== -1) {\n    lastBiblioIndex = i;\n  }\n  if (t.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (lastBiblioIndex !== -1 && copyrightIndex !== -1 && copyrightIndex > lastBiblioIndex) {\n  // Delete from the copyright paragraph back to (and including) the blank\n  // paragraph right after the bibliography's last entry, so everything in\n  // between (and including) disappears.\n  for (let i = copyrightIndex; i > lastBiblioIndex; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block (and the\n# blank paragraph immediately preceding it) that followed the last\n# bibliography entry (\"SERAFINI, Maria Jos\u00e9...\").\n#\n# The document ends with:\n#   ... SERAFINI, Maria Jos\u00e9. Como escrever textos. ...\n#   (blank paragraph)\n#   Ver no Jupiter Salvar em pdf Salvar em docx\n#   \u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ...\n#   (blank paragraph)\n#   (paragraph with page-break-before)\n#\n# After the edit the two footer paragraphs and the blank paragraph that\n# introduced them are gone, so the bibliography's last entry is followed\n# directly by the blank paragraph + page-break paragraph that used to close\n# the document.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$lastBiblioIdx = 0\n$copyrightIdx = 0\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*SERAFINI*Como escrever textos*\") { $lastBiblioIdx = $i }\n    if ($t -like \"*Contact: luizeleno@usp.br*\") { $copyrightIdx = $i }\n}\n\nif ($lastBiblioIdx -gt 0 -and $copyrightIdx -gt $lastBiblioIdx) {\n    $startRange = $d.Paragraphs.Item($lastBiblioIdx + 1).Range\n    $endRange = $d.Paragraphs.Item($copyrightIdx).Range\n    $killRange = $d.Range($startRange.Start, $endRange.End)\n    $killRange.Delete()\n}\n"}
